$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking price strings
# (e.g. "582.71") are stored as text, matching the source inlineStr cells,
# instead of being auto-coerced into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.389.17"
$ws.Range("E2").Value = "  -2.67%  "
$ws.Range("D3").Value = "2.894.39"
$ws.Range("E3").Value = "  -3.68%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "582.71"
$ws.Range("E5").Value = "  -1.82%  "
$ws.Range("D6").Value = "146.04"
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.502"
$ws.Range("E8").Value = "  -2.70%  "
$ws.Range("D9").Value = "2.891.47"
$ws.Range("E9").Value = "  -3.78%  "
$ws.Range("D10").Value = "6.60"
$ws.Range("E10").Value = "  +6.29%  "
$ws.Range("E11").Value = "  -3.52%  "
$ws.Range("E12").Value = "  -2.59%  "
$ws.Range("E13").Value = "  -3.66%  "
$ws.Range("D14").Value = "34.11"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").Value = "3.377.05"
$ws.Range("E16").Value = "  -3.48%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "60.416.19"
$ws.Range("E17").Value = "  -2.57%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "6.78"
$ws.Range("E18").Value = "  -2.89%  "
$ws.Range("D19").Value = "2.896.10"
$ws.Range("E19").Value = "  -3.57%  "
$ws.Range("D20").Value = "423.50"
$ws.Range("E20").Value = "  -5.17%  "
$ws.Range("E21").Value = "  -4.02%  "
$ws.Range("E22").Value = "  -2.66%  "
$ws.Range("E23").Value = "  -3.55%  "
$ws.Range("D24").Value = "80.59"
$ws.Range("E24").Value = "  -1.96%  "
$ws.Range("D25").Value = "10.95"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D26").Value = "2.15"
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("D27").Value = "11.78"
$ws.Range("E27").Value = "  -2.89%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "7.19"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("E31").Value = "  +4.11%  "
$ws.Range("D32").Value = "2.60"
$ws.Range("E32").Value = "  -3.52%  "
$ws.Range("D33").Value = "26.50"
$ws.Range("E33").Value = "  -3.40%  "
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("D36").Value = "1.01"
$ws.Range("E36").Value = "  -1.58%  "
$ws.Range("D37").Value = "5.63"
$ws.Range("E37").Value = "  -2.94%  "
$ws.Range("D38").Value = "49.63"
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("D40").Value = "2.93"
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("E42").Value = "  -3.06%  "
$ws.Range("D43").Value = "0.286"
$ws.Range("E43").Value = "  +2.44%  "
$ws.Range("D44").Value = "41.01"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "371.31"
$ws.Range("E45").Value = "  -5.66%  "
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("D47").Value = "133.09"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").Value = "2.643.26"
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("D50").Value = "25.06"
$ws.Range("E50").Value = "  +5.76%  "
$ws.Range("E51").Value = "  -1.05%  "

# Restore the original (default) cell style on column D now that the
# text values are committed, so no stray number-format styling remains.
$ws.Range("D2:D51").Style = "Normal"
